$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update row 2 (existing record got new values)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1553576
$ws.Range("B2").Value = 45511
$ws.Range("C2").Value = 27
$ws.Range("D2").Value = "LYY8583"
$ws.Range("E2").Value = 138.9
$ws.Range("F2").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 215415

# ---------------------------------------------------------------------------
# 2. Fill in the new rows 3-7 (same border/date style as the existing B:6 style)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 1553577
$ws.Range("B3").Value = 45509
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = "DMD2D08"
$ws.Range("E3").Value = 40.2
$ws.Range("F3").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 215250

$ws.Range("A4").Value = 1553580
$ws.Range("B4").Value = 45511
$ws.Range("C4").Value = 27
$ws.Range("D4").Value = "DMD2D08"
$ws.Range("E4").Value = 22.8
$ws.Range("F4").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 215421

$ws.Range("A5").Value = 1553581
$ws.Range("B5").Value = 45512
$ws.Range("C5").Value = 27
$ws.Range("D5").Value = "DMD2D08"
$ws.Range("E5").Value = 45.2
$ws.Range("F5").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 215632

$ws.Range("A6").Value = 1553582
$ws.Range("B6").Value = 45512
$ws.Range("C6").Value = 27
$ws.Range("D6").Value = "EGM9B58"
$ws.Range("E6").Value = 29.2
$ws.Range("F6").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 215529

$ws.Range("A7").Value = 1553583
$ws.Range("B7").Value = 45511
$ws.Range("C7").Value = 27
$ws.Range("D7").Value = "DBB6021"
$ws.Range("E7").Value = 67.8
$ws.Range("F7").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 215399

# ---------------------------------------------------------------------------
# 3. Row 8 - last row, needs a bottom border on the date cell (new style)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 1553584
$ws.Range("B8").Value = 45511
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = "HJF*0J75"
$ws.Range("E8").Value = 69.40000000000001
$ws.Range("F8").Value = "REEMBOLSO PEDAGIO"
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 215509

# give the last row's date cell its own bottom border (closing off the table)
$ws.Range("B8").Borders.Item(9).Color = 10855845
$ws.Range("B8").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 4. Column D needs an explicit width now that it has data
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 9.35

# ---------------------------------------------------------------------------
# 5. Update the selection to match the new working range
# ---------------------------------------------------------------------------
$ws.Range("A2:XFD12").Select() | Out-Null

Write-Host "done"
